# Daily attendance processing - 2025-10-08 05:44:09
# For every row in column G ("Recorded By") whose value is a comma-separated
# list that includes "System" (any case) as one of its entries, reverse the
# order of the entries in that list (e.g. "a, System" -> "System, a").
# Cells that do not contain "System" as a token are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Text

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $hasSystem = $false
    foreach ($p in $trimmed) {
        if ($p.ToLower() -eq "system") { $hasSystem = $true }
    }

    if ($hasSystem) {
        $count = $trimmed.Count
        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $trimmed[$i]
        }
        $newVal = [string]::Join(", ", $reversed)
        $cell.Value = $newVal
    }
}
